# Insert a new data row at row 310 (pushing existing rows 310:349 down to
# 311:350) on the single worksheet, then populate the new row with the
# new "Piña" price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 310..349 down to 311..350, carrying formatting along (this is
# exactly what Excel's "Insert Sheet Rows" does).
$ws.Rows.Item(310).Insert()

# Populate the newly inserted row 310 with the new record.
$ws.Range("A310").Value = 5
$ws.Range("B310").Value = "Macroferia Regional de Talca"
$ws.Range("C310").Value = "Maule"
$ws.Range("D310").Value = 44984
$ws.Range("E310").Value = 7
$ws.Range("F310").Value = "Fruta"
$ws.Range("G310").Value = 100108
$ws.Range("H310").Value = "Tropicales y subtropicales"
$ws.Range("I310").Value = 100108005
$ws.Range("J310").Value = "Piña"
$ws.Range("K310").Value = "Caramelo"
$ws.Range("L310").Value = "Segunda"
$ws.Range("M310").Value = 220
$ws.Range("N310").Value = 22000
$ws.Range("O310").Value = 22000
$ws.Range("P310").Value = 22000
$ws.Range("Q310").Value = "$/caja 14 unidades"
$ws.Range("R310").Value = "Ecuador"
$ws.Range("S310").Value = 1571
$ws.Range("T310").Value = 14
